# Auto-generated edit script: applies the numeric restatement described
# in the commit diff to the relevant Leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 767
$ws.Range("J2").Value = 1144
$ws.Range("L2").Value = 1144
$ws.Range("N2").Value = -1370
$ws.Range("H9").Value = 999359.5600000001
$ws.Range("I9").Value = 1623461.4
$ws.Range("K9").Value = 1623461.4
$ws.Range("M9").Value = -1623292.4
$ws.Range("H18").Value = 2059.8
$ws.Range("I18").Value = 800
$ws.Range("K18").Value = 800
$ws.Range("M18").Value = -516
$ws.Range("H32").Value = 14886.223
$ws.Range("I32").Value = 17997.4
$ws.Range("K32").Value = 17997.4
$ws.Range("M32").Value = -17671.4
$ws.Range("H40").Value = 11585
$ws.Range("I40").Value = 6981.25
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 6981.25
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = -6806.25
$ws.Range("N40").Value = -30350
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H70").Value = 174566.5
$ws.Range("I70").Value = 4133.3335
$ws.Range("K70").Value = 12400.0005
$ws.Range("M70").Value = -12130.0005
$ws.Range("H73").Value = 174566.5
$ws.Range("I73").Value = 4133.3335
$ws.Range("K73").Value = 12400.0005
$ws.Range("M73").Value = -11464.0005
$ws.Range("H118").Value = 924.75
$ws.Range("I118").Value = 924.75
$ws.Range("K118").Value = 2774.25
$ws.Range("M118").Value = -1117.25
$ws.Range("H132").Value = 3800.3572
$ws.Range("I132").Value = 3785.3333
$ws.Range("J132").Value = 4206
$ws.Range("K132").Value = 11355.9999
$ws.Range("L132").Value = 12618
$ws.Range("M132").Value = -8825.999899999999
$ws.Range("N132").Value = -17678
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3560.575
$ws.Range("I32").Value = 1606.0541
$ws.Range("K32").Value = 1606.0541
$ws.Range("M32").Value = -1319.0541
$ws.Range("H74").Value = 3285.7144
$ws.Range("I74").Value = 3285.7144
$ws.Range("K74").Value = 3285.7144
$ws.Range("M74").Value = -2411.7144
$ws.Range("H77").Value = 3285.7144
$ws.Range("I77").Value = 3285.7144
$ws.Range("K77").Value = 16428.572
$ws.Range("M77").Value = -12060.572

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 547.55554
$ws.Range("I22").Value = 571.8333
$ws.Range("K22").Value = 571.8333
$ws.Range("M22").Value = -398.8333
$ws.Range("H86").Value = 3750.5
$ws.Range("I86").Value = 1579.6
$ws.Range("K86").Value = 1579.6
$ws.Range("M86").Value = -456.5999999999999
$ws.Range("H89").Value = 3750.5
$ws.Range("I89").Value = 1579.6
$ws.Range("K89").Value = 7898
$ws.Range("M89").Value = -2282

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1170.8
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 6924
$ws.Range("I31").Value = 8299.777
$ws.Range("J31").Value = 3828.5
$ws.Range("K31").Value = 8299.777
$ws.Range("L31").Value = 3828.5
$ws.Range("M31").Value = -8004.777
$ws.Range("N31").Value = -4418.5
$ws.Range("H34").Value = 6924
$ws.Range("I34").Value = 8299.777
$ws.Range("J34").Value = 3828.5
$ws.Range("K34").Value = 8299.777
$ws.Range("L34").Value = 3828.5
$ws.Range("M34").Value = -8097.777
$ws.Range("N34").Value = -4232.5
$ws.Range("H92").Value = 44833.332
$ws.Range("J92").Value = 44833.332
$ws.Range("L92").Value = 44833.332
$ws.Range("N92").Value = -49825.332
$ws.Range("H107").Value = 1095.2667
$ws.Range("I107").Value = 994.4167
$ws.Range("K107").Value = 994.4167
$ws.Range("M107").Value = 925.5833
$ws.Range("H132").Value = 5672.479
$ws.Range("I132").Value = 4501.567
$ws.Range("J132").Value = 7624
$ws.Range("K132").Value = 13504.701
$ws.Range("L132").Value = 22872
$ws.Range("M132").Value = -10974.701
$ws.Range("N132").Value = -27932

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1808.4
$ws.Range("I34").Value = 102.22222
$ws.Range("J34").Value = 2768.125
$ws.Range("K34").Value = 306.66666
$ws.Range("L34").Value = 8304.375
$ws.Range("M34").Value = -222.66666
$ws.Range("N34").Value = -8472.375
$ws.Range("H39").Value = 8253
$ws.Range("I39").Value = 2147.5
$ws.Range("J39").Value = 9997.429
$ws.Range("K39").Value = 6442.5
$ws.Range("L39").Value = 29992.287
$ws.Range("M39").Value = -6148.5
$ws.Range("N39").Value = -30580.287
$ws.Range("H64").Value = 12197
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 14746.25
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 44238.75
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -44778.75
$ws.Range("H67").Value = 12197
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 14746.25
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 44238.75
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -46110.75
$ws.Range("H80").Value = 3960
$ws.Range("J80").Value = 4200
$ws.Range("L80").Value = 12600
$ws.Range("N80").Value = -14472
$ws.Range("H83").Value = 3960
$ws.Range("J83").Value = 4200
$ws.Range("L83").Value = 37800
$ws.Range("N83").Value = -47160
$ws.Range("H107").Value = 903.06665
$ws.Range("I107").Value = 810.44446
$ws.Range("J107").Value = 1042
$ws.Range("K107").Value = 2431.33338
$ws.Range("L107").Value = 3126
$ws.Range("M107").Value = -511.33338
$ws.Range("N107").Value = -6966
$ws.Range("H114").Value = 1425.4
$ws.Range("J114").Value = 1524.75
$ws.Range("L114").Value = 4574.25
$ws.Range("N114").Value = -11082.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H122").Value = 3860.5715
$ws.Range("I122").Value = 3767.6
$ws.Range("J122").Value = 4093
$ws.Range("K122").Value = 11302.8
$ws.Range("L122").Value = 12279
$ws.Range("M122").Value = -8852.799999999999
$ws.Range("N122").Value = -17179

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224
$ws.Range("H16").Value = 834.9231
$ws.Range("I16").Value = 834.9231
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 834.9231
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -664.9231
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 4285.2856
$ws.Range("I22").Value = 3999.5
$ws.Range("J22").Value = 4399.6
$ws.Range("K22").Value = 3999.5
$ws.Range("L22").Value = 4399.6
$ws.Range("M22").Value = -3704.5
$ws.Range("N22").Value = -4989.6
$ws.Range("H27").Value = 4285.2856
$ws.Range("I27").Value = 3999.5
$ws.Range("J27").Value = 4399.6
$ws.Range("K27").Value = 3999.5
$ws.Range("L27").Value = 4399.6
$ws.Range("M27").Value = -3892.5
$ws.Range("N27").Value = -4613.6
$ws.Range("H46").Value = 4590.273
$ws.Range("I46").Value = 5638.2856
$ws.Range("K46").Value = 5638.2856
$ws.Range("M46").Value = -5450.2856
$ws.Range("H68").Value = 6997.6
$ws.Range("I68").Value = 4993
$ws.Range("K68").Value = 4993
$ws.Range("M68").Value = -4244
$ws.Range("H71").Value = 6997.6
$ws.Range("I71").Value = 4993
$ws.Range("K71").Value = 24965
$ws.Range("M71").Value = -21221
$ws.Range("H93").Value = 1928.5555
$ws.Range("I93").Value = 2044.125
$ws.Range("J93").Value = 1004
$ws.Range("K93").Value = 2044.125
$ws.Range("L93").Value = 1004
$ws.Range("M93").Value = -796.125
$ws.Range("N93").Value = -3500
$ws.Range("H100").Value = 5326.3335
$ws.Range("I100").Value = 4599
$ws.Range("K100").Value = 4599
$ws.Range("M100").Value = -4058
$ws.Range("H133").Value = 94999
$ws.Range("J133").Value = 94999
$ws.Range("L133").Value = 94999
$ws.Range("N133").Value = -100059
$ws.Range("H136").Value = 1696.9642
$ws.Range("I136").Value = 1543.579
$ws.Range("J136").Value = 2020.7778
$ws.Range("K136").Value = 4630.737
$ws.Range("L136").Value = 6062.3334
$ws.Range("M136").Value = -2080.737
$ws.Range("N136").Value = -11162.3334

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5285.5713
$ws.Range("J2").Value = 5285.5713
$ws.Range("L2").Value = 5285.5713
$ws.Range("N2").Value = -5509.5713
$ws.Range("H3").Value = 3337333.2
$ws.Range("I3").Value = 3337333.2
$ws.Range("K3").Value = 3337333.2
$ws.Range("M3").Value = -3337219.2
$ws.Range("H122").Value = 3285.0908
$ws.Range("I122").Value = 3169.7144
$ws.Range("K122").Value = 9509.143199999999
$ws.Range("M122").Value = -7059.143199999999
$ws.Range("H132").Value = 3445.4614
$ws.Range("I132").Value = 3708.2727
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 11124.8181
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -8594.8181
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 3238.375
$ws.Range("I136").Value = 3164.842
$ws.Range("J136").Value = 3517.8
$ws.Range("K136").Value = 9494.526
$ws.Range("L136").Value = 10553.4
$ws.Range("M136").Value = -6944.526
$ws.Range("N136").Value = -15653.4
